$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) updates
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 93
$wsExpo.Range("F4").Value = 267
$wsExpo.Range("F6").Value = 39
$wsExpo.Range("F7").Value = 264
$wsExpo.Range("F9").Value = 1960
$wsExpo.Range("F11").Value = 4621
$wsExpo.Range("F12").Value = 79
$wsExpo.Range("F13").Value = 325

# Sheet "全部类型" (all types) updates - same events, different row numbers
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 93
$wsAll.Range("F6").Value = 267
$wsAll.Range("F8").Value = 39
$wsAll.Range("F9").Value = 264
$wsAll.Range("F13").Value = 1960
$wsAll.Range("F15").Value = 4621
$wsAll.Range("F16").Value = 79
$wsAll.Range("F17").Value = 325
